$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for a new weekly observation: insert a blank row at position 3
# (this leaves row 2's current data untouched and only copies the plain
# date-column style down, avoiding the header row's formatting).
$ws.Rows.Item(3).Insert()

# Push the data that used to live in row 2 down into the newly created row 3.
$ws.Range("A2:R2").Copy($ws.Range("A3:R3"))

# Overwrite row 2 with the latest weekly price observation.
$ws.Range("A2").Value = 6
$ws.Range("B2").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C2").Value = "Metropolitana"
$ws.Range("D2").Value = 44599
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = 100114007
$ws.Range("G2").Value = "Jengibre"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 13000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 14200
$ws.Range("N2").Value = "$/caja 13 kilos"
$ws.Range("O2").Value = "Perú"
$ws.Range("P2").Value = 1092
$ws.Range("Q2").Value = 13
$ws.Range("R2").Value = "Hortaliza"
